# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation" on all sheets
# - Shrink the (now shorter) Status columns to match the new text width

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

# Update the localization status values
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Re-fit the status columns now that the text is shorter
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
